$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "0, 0"
$ws.Range("D2").Value = "4, 3"
$ws.Range("D4").Value = "4, 16"

$range = $ws.Range("A1:F4")
$range.Borders.LineStyle = 1
$range.HorizontalAlignment = -4108

$ws.Range("D4").Select()
